$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AH column changes (value + fill color / style) ---
$ws.Range("AH18").Interior.Color = 16777215
$ws.Range("AH18").Value = 20
$ws.Range("AH20").Value = 20
$ws.Range("AH21").Interior.Color = 32768
$ws.Range("AH21").Value = 40
$ws.Range("AH22").Value = 20
$ws.Range("AH24").Value = 32
$ws.Range("AH25").Value = 30
$ws.Range("AH29").Interior.Color = 65535
$ws.Range("AH29").Value = 1
$ws.Range("AH30").Interior.Color = 255
$ws.Range("AH30").Value = 0
$ws.Range("AH31").Interior.Color = 32768
$ws.Range("AH31").Value = 34
$ws.Range("AH32").Interior.Color = 255
$ws.Range("AH32").Value = 0
$ws.Range("AH36").Interior.Color = 255
$ws.Range("AH36").Value = 0
$ws.Range("AH38").Value = 30
$ws.Range("AH39").Value = 30
$ws.Range("AH41").Interior.Color = 16777215
$ws.Range("AH41").Value = 20
$ws.Range("AH45").Value = 21
$ws.Range("AH47").Interior.Color = 16777215
$ws.Range("AH47").Value = 30
$ws.Range("AH48").Value = 30
$ws.Range("AH49").Interior.Color = 16777215
$ws.Range("AH49").Value = 30
$ws.Range("AH51").Interior.Color = 255
$ws.Range("AH51").Value = 0
$ws.Range("AH53").Value = 8
$ws.Range("AH54").Interior.Color = 255
$ws.Range("AH54").Value = 0
$ws.Range("AH55").Interior.Color = 65535
$ws.Range("AH55").Value = 1
$ws.Range("AH56").Value = 30
$ws.Range("AH58").Interior.Color = 255
$ws.Range("AH58").Value = 0
$ws.Range("AH60").Value = 25
$ws.Range("AH63").Interior.Color = 16777215
$ws.Range("AH63").Value = 20
$ws.Range("AH64").Interior.Color = 16777215
$ws.Range("AH64").Value = 20
$ws.Range("AH65").Value = 23
$ws.Range("AH69").Interior.Color = 255
$ws.Range("AH69").Value = 0
$ws.Range("AH76").Value = 5
$ws.Range("AH118").Value = 20

# --- AI column changes (value only, keep text/number type as in target) ---
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "2500"
$ws.Range("AI5").NumberFormat = "@"
$ws.Range("AI5").Value = "0"
$ws.Range("AI6").NumberFormat = "@"
$ws.Range("AI6").Value = "0"
$ws.Range("AI18").NumberFormat = "@"
$ws.Range("AI18").Value = "2780"
$ws.Range("AI19").NumberFormat = "@"
$ws.Range("AI19").Value = "2498"
$ws.Range("AI20").NumberFormat = "@"
$ws.Range("AI20").Value = "2783"
$ws.Range("AI21").NumberFormat = "@"
$ws.Range("AI21").Value = "2975"
$ws.Range("AI22").NumberFormat = "@"
$ws.Range("AI22").Value = "2814"
$ws.Range("AI23").NumberFormat = "@"
$ws.Range("AI23").Value = "2870"
$ws.Range("AI24").NumberFormat = "@"
$ws.Range("AI24").Value = "2938"
$ws.Range("AI25").NumberFormat = "@"
$ws.Range("AI25").Value = "2952"
$ws.Range("AI27").NumberFormat = "@"
$ws.Range("AI27").Value = "0"
$ws.Range("AI29").NumberFormat = "@"
$ws.Range("AI29").Value = "2516"
$ws.Range("AI30").NumberFormat = "@"
$ws.Range("AI30").Value = "2500"
$ws.Range("AI31").NumberFormat = "@"
$ws.Range("AI31").Value = "2995"
$ws.Range("AI32").NumberFormat = "@"
$ws.Range("AI32").Value = "2498"
$ws.Range("AI35").Value = 5280
$ws.Range("AI36").NumberFormat = "@"
$ws.Range("AI36").Value = "2499"
$ws.Range("AI37").NumberFormat = "@"
$ws.Range("AI37").Value = "2950"
$ws.Range("AI38").NumberFormat = "@"
$ws.Range("AI38").Value = "2985"
$ws.Range("AI39").NumberFormat = "@"
$ws.Range("AI39").Value = "2890"
$ws.Range("AI41").NumberFormat = "@"
$ws.Range("AI41").Value = "2927"
$ws.Range("AI42").NumberFormat = "@"
$ws.Range("AI42").Value = "2498"
$ws.Range("AI45").NumberFormat = "@"
$ws.Range("AI45").Value = "2667"
$ws.Range("AI46").NumberFormat = "@"
$ws.Range("AI46").Value = "2740"
$ws.Range("AI47").NumberFormat = "@"
$ws.Range("AI47").Value = "3005"
$ws.Range("AI48").NumberFormat = "@"
$ws.Range("AI48").Value = "3015"
$ws.Range("AI49").NumberFormat = "@"
$ws.Range("AI49").Value = "3048"
$ws.Range("AI50").NumberFormat = "@"
$ws.Range("AI50").Value = "2970"
$ws.Range("AI51").NumberFormat = "@"
$ws.Range("AI51").Value = "2573"
$ws.Range("AI52").NumberFormat = "@"
$ws.Range("AI52").Value = "3045"
$ws.Range("AI53").NumberFormat = "@"
$ws.Range("AI53").Value = "2689"
$ws.Range("AI54").NumberFormat = "@"
$ws.Range("AI54").Value = "2518"
$ws.Range("AI55").NumberFormat = "@"
$ws.Range("AI55").Value = "2529"
$ws.Range("AI56").NumberFormat = "@"
$ws.Range("AI56").Value = "3091"
$ws.Range("AI57").NumberFormat = "@"
$ws.Range("AI57").Value = "2732"
$ws.Range("AI58").NumberFormat = "@"
$ws.Range("AI58").Value = "2587"
$ws.Range("AI59").NumberFormat = "@"
$ws.Range("AI59").Value = "2765"
$ws.Range("AI60").NumberFormat = "@"
$ws.Range("AI60").Value = "2927"
$ws.Range("AI62").NumberFormat = "@"
$ws.Range("AI62").Value = "2752"
$ws.Range("AI63").NumberFormat = "@"
$ws.Range("AI63").Value = "2787"
$ws.Range("AI64").NumberFormat = "@"
$ws.Range("AI64").Value = "2866"
$ws.Range("AI65").NumberFormat = "@"
$ws.Range("AI65").Value = "2729"
$ws.Range("AI68").NumberFormat = "@"
$ws.Range("AI68").Value = "0"
$ws.Range("AI69").NumberFormat = "@"
$ws.Range("AI69").Value = "2499"
$ws.Range("AI70").NumberFormat = "@"
$ws.Range("AI70").Value = "0"
$ws.Range("AI73").NumberFormat = "@"
$ws.Range("AI73").Value = "0"
$ws.Range("AI76").NumberFormat = "@"
$ws.Range("AI76").Value = "2594"
$ws.Range("AI77").NumberFormat = "@"
$ws.Range("AI77").Value = "2519"
$ws.Range("AI78").NumberFormat = "@"
$ws.Range("AI78").Value = "0"
$ws.Range("AI81").NumberFormat = "@"
$ws.Range("AI81").Value = "0"
$ws.Range("AI83").NumberFormat = "@"
$ws.Range("AI83").Value = "2514"
$ws.Range("AI84").NumberFormat = "@"
$ws.Range("AI84").Value = "0"
$ws.Range("AI90").NumberFormat = "@"
$ws.Range("AI90").Value = "1000"
$ws.Range("AI96").NumberFormat = "@"
$ws.Range("AI96").Value = "0"
$ws.Range("AI118").NumberFormat = "@"
$ws.Range("AI118").Value = "2858"
$ws.Range("AI119").NumberFormat = "@"
$ws.Range("AI119").Value = "0"
$ws.Range("AI120").NumberFormat = "@"
$ws.Range("AI120").Value = "0"
